$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - new TPM-derived values
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 1.369418292493333
$ws.Range("R2").Value = 12.32476463244
$ws.Range("S2").Value = 0.2403816673726824
$ws.Range("T2").Value = 0.2403816673726824

# Row 3 - recomputed specificity values
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("S3").Value = 0.3167483425780597
$ws.Range("T3").Value = 0.3167483425780597

# Row 4 - recomputed specificity values
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("S4").Value = 0.4428699900492579
$ws.Range("T4").Value = 0.4428699900492579
